$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20. This shifts the existing rows 20-31
# down to 21-32, preserving all of their data (dates, prices, origin, etc.)
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new record
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44755
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112040
$ws.Range("G20").Value = "Cilantro"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 8500
$ws.Range("M20").Value = 8500
$ws.Range("N20").Value = "`$/caja 36 atados"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 236
$ws.Range("Q20").Value = 36
$ws.Range("R20").Value = "Hortaliza"
